$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 23; everything below shifts down by one.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new record.
$ws.Cells.Item(23, 1).Value = 6
$ws.Cells.Item(23, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(23, 3).Value = "Metropolitana"
$ws.Cells.Item(23, 4).Value = 45133
$ws.Cells.Item(23, 5).Value = 13
$ws.Cells.Item(23, 6).Value = 100112035
$ws.Cells.Item(23, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 310
$ws.Cells.Item(23, 11).Value = 18000
$ws.Cells.Item(23, 12).Value = 20000
$ws.Cells.Item(23, 13).Value = 19032
$ws.Cells.Item(23, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(23, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(23, 16).Value = 1269
$ws.Cells.Item(23, 17).Value = 15
$ws.Cells.Item(23, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date number format used by the rest of column D.
$ws.Cells.Item(23, 4).NumberFormat = $ws.Cells.Item(24, 4).NumberFormat
